$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" data column (Q), matching the
# formatting already used by the adjacent "2021" column (P).
$ws.Range("P4:P5").Copy()
$ws.Range("Q4:Q5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Q4").Value = 2022
$ws.Range("Q5").Value = 64.2

$excel.CutCopyMode = 0

# Move the active selection, as recorded in the saved workbook view.
$ws.Range("R4").Select()
